# Auto-generated edit script: updates leve-profit calculation cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per the refreshed
# Universalis pricing snapshot pulled by the scheduled runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2440.4
$ws.Range("I6").Value = 3466.6667
$ws.Range("J6").Value = 901
$ws.Range("K6").Value = 10400.0001
$ws.Range("L6").Value = 2703
$ws.Range("M6").Value = -10288.0001
$ws.Range("N6").Value = -2927

$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()

$ws.Range("H132").Value = 228923.19
$ws.Range("I132").Value = 1884.9354
$ws.Range("J132").Value = 770322.0600000001
$ws.Range("K132").Value = 5654.8062
$ws.Range("L132").Value = 2310966.18
$ws.Range("M132").Value = -3124.8062
$ws.Range("N132").Value = -2316026.18

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 3050
$ws.Range("I88").Value = 4000
$ws.Range("J88").Value = 2733.3333
$ws.Range("K88").Value = 4000
$ws.Range("L88").Value = 2733.3333
$ws.Range("M88").Value = -3594
$ws.Range("N88").Value = -3545.3333

$ws.Range("H91").Value = 3050
$ws.Range("I91").Value = 4000
$ws.Range("J91").Value = 2733.3333
$ws.Range("K91").Value = 4000
$ws.Range("L91").Value = 2733.3333
$ws.Range("M91").Value = -2596
$ws.Range("N91").Value = -5541.3333

$ws.Range("H110").Value = 953.6667
$ws.Range("I110").Value = 953.6667
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 953.6667
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = 1091.3333
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 1700
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1700
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1700
$ws.Range("N5").Value = -1926
$ws.Range("M5").ClearContents()

$ws.Range("H86").Value = 2245.1538
$ws.Range("I86").Value = 1708.8889
$ws.Range("J86").Value = 3451.75
$ws.Range("K86").Value = 1708.8889
$ws.Range("L86").Value = 3451.75
$ws.Range("M86").Value = -585.8888999999999
$ws.Range("N86").Value = -5697.75

$ws.Range("H89").Value = 2245.1538
$ws.Range("I89").Value = 1708.8889
$ws.Range("J89").Value = 3451.75
$ws.Range("K89").Value = 8544.4445
$ws.Range("L89").Value = 17258.75
$ws.Range("M89").Value = -2928.4445
$ws.Range("N89").Value = -28490.75

$ws.Range("H102").Value = 15000
$ws.Range("I102").Value = 15000
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 15000
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -11755
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 23339.166
$ws.Range("I107").Value = 1674
$ws.Range("J107").Value = 45004.332
$ws.Range("K107").Value = 1674
$ws.Range("L107").Value = 45004.332
$ws.Range("M107").Value = 246
$ws.Range("N107").Value = -48844.332

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H17").Value = 5000
$ws.Range("J17").Value = 5000
$ws.Range("L17").Value = 5000
$ws.Range("N17").Value = -5348

$ws.Range("H99").Value = 1875.6428
$ws.Range("I99").Value = 1526
$ws.Range("J99").Value = 2389.8235
$ws.Range("K99").Value = 1526
$ws.Range("L99").Value = 2389.8235
$ws.Range("M99").Value = -28
$ws.Range("N99").Value = -5385.8235

$ws.Range("H126").Value = 1875.6428
$ws.Range("I126").Value = 1526
$ws.Range("J126").Value = 2389.8235
$ws.Range("K126").Value = 4578
$ws.Range("L126").Value = 7169.470499999999
$ws.Range("M126").Value = -2108
$ws.Range("N126").Value = -12109.4705

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 321.25
$ws.Range("I7").Value = 196
$ws.Range("J7").Value = 530
$ws.Range("K7").Value = 588
$ws.Range("L7").Value = 1590
$ws.Range("M7").Value = -476
$ws.Range("N7").Value = -1814

$ws.Range("H80").Value = 7299.6
$ws.Range("I80").Value = 3500
$ws.Range("J80").Value = 8249.5
$ws.Range("K80").Value = 10500
$ws.Range("L80").Value = 24748.5
$ws.Range("M80").Value = -9564
$ws.Range("N80").Value = -26620.5

$ws.Range("H83").Value = 7299.6
$ws.Range("I83").Value = 3500
$ws.Range("J83").Value = 8249.5
$ws.Range("K83").Value = 31500
$ws.Range("L83").Value = 74245.5
$ws.Range("M83").Value = -26820
$ws.Range("N83").Value = -83605.5

$ws.Range("H92").Value = 625
$ws.Range("I92").Value = 500
$ws.Range("J92").Value = 750
$ws.Range("K92").Value = 1500
$ws.Range("L92").Value = 2250
$ws.Range("M92").Value = -252
$ws.Range("N92").Value = -4746

$ws.Range("H110").Value = 10114.111
$ws.Range("I110").Value = 342.33334
$ws.Range("K110").Value = 1027.00002
$ws.Range("M110").Value = 3062.99998

$ws.Range("H112").Value = 3283.7
$ws.Range("I112").Value = 902.3333
$ws.Range("J112").Value = 4304.2856
$ws.Range("K112").Value = 2706.9999
$ws.Range("L112").Value = 12912.8568
$ws.Range("M112").Value = -1598.9999
$ws.Range("N112").Value = -15128.8568

$ws.Range("H117").Value = 1980.85
$ws.Range("I117").Value = 1379
$ws.Range("J117").Value = 2047.7222
$ws.Range("K117").Value = 4137
$ws.Range("L117").Value = 6143.1666
$ws.Range("M117").Value = -695
$ws.Range("N117").Value = -13027.1666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 900
$ws.Range("I102").Value = 900
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 900
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 722
$ws.Range("N102").ClearContents()

$ws.Range("H107").Value = 1177031.5
$ws.Range("I107").Value = 2941808.2
$ws.Range("K107").Value = 2941808.2
$ws.Range("M107").Value = -2939888.2

$ws.Range("H126").Value = 3305
$ws.Range("I126").Value = 2103
$ws.Range("J126").Value = 4507
$ws.Range("K126").Value = 6309
$ws.Range("L126").Value = 13521
$ws.Range("M126").Value = -3839
$ws.Range("N126").Value = -18461

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2312.625
$ws.Range("I68").Value = 2200.2856
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 2200.2856
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = -1451.2856
$ws.Range("N68").Value = -3898

$ws.Range("H71").Value = 2312.625
$ws.Range("I71").Value = 2200.2856
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 11001.428
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = -7257.428
$ws.Range("N71").Value = -19488

$ws.Range("H122").Value = 4499.4546
$ws.Range("J122").Value = 3003.3333
$ws.Range("L122").Value = 9009.999899999999
$ws.Range("N122").Value = -13909.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1774.1333
$ws.Range("J81").Value = 1898.8572
$ws.Range("L81").Value = 3797.7144
$ws.Range("N81").Value = -5919.7144

$ws.Range("H84").Value = 1774.1333
$ws.Range("J84").Value = 1898.8572
$ws.Range("L84").Value = 18988.572
$ws.Range("N84").Value = -29596.572

$ws.Range("H126").Value = 674.25
$ws.Range("I126").Value = 644.63635
$ws.Range("J126").Value = 1000
$ws.Range("K126").Value = 1933.90905
$ws.Range("L126").Value = 3000
$ws.Range("M126").Value = 536.09095
$ws.Range("N126").Value = -7940

$ws.Range("H132").Value = 1645.0526
$ws.Range("I132").Value = 1211.7693
$ws.Range("J132").Value = 2583.8333
$ws.Range("K132").Value = 3635.3079
$ws.Range("L132").Value = 7751.499899999999
$ws.Range("M132").Value = -1105.3079
$ws.Range("N132").Value = -12811.4999

